# Update project assignment sheet: mark the next batch of Apache Commons
# projects (rows 58-78, column C "status") as "DONE" instead of "IN PROGRESS".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C58:C78").Value = "DONE"

# Reflect the author's updated scroll position / active selection in the sheet.
$excel.ActiveWindow.ScrollRow = 54
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("K69").Select()
